# "Add DH Storage and tidy some files"
#
# On the "multi" sheet, two new separator rows (the same shape as the
# existing storage-technology separator rows, e.g. row 8 / row 10 - a lone
# J cell holding -1 plus a blank, s="17"-styled K cell) are inserted, one
# above the RSD 2021-2025 budget row and one above the SRV 2021-2025 budget
# row, pushing the remaining sector-budget rows down. The "multi" tab also
# becomes the active tab/sheet (replacing "config" as the selected tab).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("multi")

# --- Insert separator row above what is currently row 12 (RSD 2021-2025) ---
$ws.Range("A12:L12").EntireRow.Insert()
$ws.Range("J12").Value = -1

# --- Insert separator row above what is currently row 13 (SRV 2021-2025),
#     i.e. row 14 after the previous insert shifted everything down one ---
$ws.Range("A14:L14").EntireRow.Insert()
$ws.Range("J14").Value = -1

# Update the sheet's remembered selection to match the new layout.
$ws.Range("J15").Select()

# Make "multi" the active/selected sheet (was "config").
$ws.Select()
